$d = $word.ActiveDocument

$replacements = @(
    @{old = "310×6=1860"; new = "106×9=954"},
    @{old = "311×6=1866"; new = "848×9=7632"},
    @{old = "176×8=1408"; new = "901×7=6307"},
    @{old = "288×7=2016"; new = "955×8=7640"},
    @{old = "222×5=1110"; new = "515×9=4635"},
    @{old = "257×2=514"; new = "734×7=5138"},
    @{old = "716×7=5012"; new = "872×4=3488"},
    @{old = "740×5=3700"; new = "428×3=1284"},
    @{old = "873×3=2619"; new = "505×6=3030"},
    @{old = "392×7=2744"; new = "428×2=856"},
    @{old = "493×9=4437"; new = "423×8=3384"},
    @{old = "150×9=1350"; new = "985×8=7880"},
    @{old = "509×6=3054"; new = "480×3=1440"},
    @{old = "229×2=458"; new = "286×2=572"},
    @{old = "598×4=2392"; new = "876×3=2628"},
    @{old = "772×4=3088"; new = "571×4=2284"},
    @{old = "288×5=1440"; new = "534×4=2136"},
    @{old = "845×5=4225"; new = "539×7=3773"},
    @{old = "266×4=1064"; new = "993×2=1986"},
    @{old = "538×8=4304"; new = "465×4=1860"},
    @{old = "410×4=1640"; new = "418×7=2926"},
    @{old = "894×7=6258"; new = "987×6=5922"},
    @{old = "227×2=454"; new = "296×8=2368"},
    @{old = "676×6=4056"; new = "103×3=309"},
    @{old = "141×8=1128"; new = "909×6=5454"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}

Write-Output "Done"
